# Fixed a bug in genreels
# The rows of reel-weight data (rows 3-21, columns A-F) were being
# written in the wrong order. This re-applies the corrected ordering
# produced by the fixed genreels routine.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values (symbol, reel1, reel2, reel3, reel4, reel5) for each row,
# after the genreels ordering fix. Row 2 and rows 22-25 are unaffected.
$data = @{
    3  = @(901, 16, 15, 45, 60, 60)
    4  = @(401, 9, 48, 67, 75, 45)
    5  = @(601, 9, 60, 67, 60, 42)
    6  = @(501, 9, 52, 30, 75, 45)
    7  = @(1201, 2, 10, 10, 10, 10)
    8  = @(1203, 3, 15, 15, 15, 15)
    9  = @(1001, 18, 30, 75, 60, 72)
    10 = @(701, 3, 90, 45, 97, 15)
    11 = @(201, 9, 30, 15, 45, 30)
    12 = @(1202, 2, 10, 10, 10, 10)
    13 = @(101, 9, 30, 15, 60, 15)
    14 = @(902, 1, 0, 0, 0, 0)
    15 = @(301, 6, 45, 30, 60, 45)
    16 = @(1, 0, 2, 2, 2, 2)
    17 = @(3, 0, 3, 3, 3, 3)
    18 = @(502, 0, 4, 0, 0, 0)
    19 = @(1101, 0, 15, 30, 30, 0)
    20 = @(802, 0, 4, 5, 4, 0)
    21 = @(2, 0, 2, 2, 2, 2)
}

foreach ($rowNum in $data.Keys) {
    $values = $data[$rowNum]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($rowNum, $col).Value = $values[$col - 1]
    }
}
